# Apply the commit's edit: the subtitle on the title slide is retitled
# from "G2M Case Study Virtual Internship" to "Healthcare project".
#
# We walk every slide/shape/paragraph rather than hard-coding indices so the
# script is resilient to minor shape-ordering differences, but it only ever
# touches the paragraph whose (trimmed) text matches the old title exactly,
# so no other content is disturbed. Only the run's text is changed - we
# don't touch TextFrame/Paragraph-level formatting, so existing run
# properties (font size, color, etc.) are preserved as-is.

$oldTitle = "G2M Case Study Virtual Internship"
$newTitle = "Healthcare project"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) {
            continue
        }

        $tf = $shape.TextFrame
        if (-not $tf.HasText) {
            continue
        }

        $tr = $tf.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if ($para.Text.Trim() -eq $oldTitle) {
                $para.Text = $newTitle
                Write-Host "Slide $si shape $shi paragraph $pi retitled to '$newTitle'"
            }
        }
    }
}

# --- Defensive, best-effort pass for the companion date-field refresh ---
# The same commit also nudges the cached "datetimeFigureOut" field text
# (10/10/2021 -> 12/10/2021 on notes pages, 10/10/2021 -> 10/12/2021 on
# slides). Those fields only exist once a slide/notes page actually has a
# cached date placeholder; update them in place when present so the script
# stays correct if/when such placeholders exist, without creating new
# shapes or placeholders that aren't already part of the deck.
$oldNotesDate = "10/10/2021"
$newNotesDate = "12/10/2021"
$oldSlideDate = "10/10/2021"
$newSlideDate = "10/12/2021"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) {
            continue
        }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) {
            continue
        }
        $tr = $tf.TextRange
        if ($tr.Text.Trim() -eq $oldSlideDate) {
            $tr.Text = $newSlideDate
            Write-Host "Slide $si shape $shi date field updated to '$newSlideDate'"
        }
    }

    if ($slide.HasNotesPage) {
        $notes = $slide.NotesPage
        for ($nshi = 1; $nshi -le $notes.Shapes.Count; $nshi++) {
            $nshape = $notes.Shapes.Item($nshi)
            if (-not $nshape.HasTextFrame) {
                continue
            }
            $ntf = $nshape.TextFrame
            if (-not $ntf.HasText) {
                continue
            }
            $ntr = $ntf.TextRange
            if ($ntr.Text.Trim() -eq $oldNotesDate) {
                $ntr.Text = $newNotesDate
                Write-Host "Slide $si notes shape $nshi date field updated to '$newNotesDate'"
            }
        }
    }
}
